$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# Version 5.0.0 -> 6.0.0
$wsMeta.Range("B3").Value = "6.0.0"

# Date changed
$wsMeta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$wsMeta.Range("B9").Value = "Alvearie Team"

# Contact row becomes Jurisdiction / United States of America
$wsMeta.Range("A10").Value = "Jurisdiction"
$wsMeta.Range("B10").Value = "United States of America"

# Remove the duplicate "Contact" row (old row 11)
$wsMeta.Rows.Item(11).Delete()

# Elements sheet: root element Short/Definition text
$wsElem.Range("K2").Value = "Source Event Timestamp"
$wsElem.Range("L2").Value = "Date and time of the source event that triggers either the creation or updating of this FHIR resource"
